# Lesson 24 grammar slide: bump the "PRESENT / PAST / FUTURE" tense labels
# and all "work / worked / will work" cells from 24pt to 28pt, and resize
# the (auto-fit) textboxes to the new on-screen box PowerPoint computes for
# the larger font (values captured from the canonical, already-resized XML).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Resize-Label {
    param(
        [string]$Name,
        [double]$Left,
        [double]$Top,
        [double]$Width,
        [double]$Height
    )

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq $Name) {
            $shp.TextFrame.TextRange.Font.Size = 28
            $shp.Left = $Left
            $shp.Top = $Top
            $shp.Width = $Width
            $shp.Height = $Height
            break
        }
    }
}

# Row 1 - PRESENT / PAST / FUTURE headers
Resize-Label "TextBox 54" 287.09625 119.09484 127.483   33.92815
Resize-Label "TextBox 55" 481.4061  119.09484 73.713    33.92815
Resize-Label "TextBox 67" 628.95123 119.09484 113.4725  33.92815

# Row 2 (y=2266321 EMU)
Resize-Label "TextBox 8"  317.89406 176.0271  65.8874   33.92815
Resize-Label "TextBox 10" 468.7841  176.0271  98.9571   33.92815
Resize-Label "TextBox 11" 624.15485 176.0271  123.0652  33.92815

# Row 3 (y=2976142 EMU)
Resize-Label "TextBox 14" 317.89406 231.91855 65.8874   33.92815
Resize-Label "TextBox 15" 468.7841  231.91855 98.9571   33.92815
Resize-Label "TextBox 16" 624.15485 231.91855 123.0652  33.92815

# Row 4 (y=3712400 EMU)
Resize-Label "TextBox 17" 310.9519  289.8916  79.7716   33.92815
Resize-Label "TextBox 18" 468.7841  289.8916  98.9571   33.92815
Resize-Label "TextBox 19" 624.15485 289.8916  123.0652  33.92815

# Row 5 (y=4432241 EMU)
Resize-Label "TextBox 20" 317.89406 346.572   65.8874   33.92815
Resize-Label "TextBox 21" 468.7841  346.572   98.9571   33.92815
Resize-Label "TextBox 22" 624.15485 346.572   123.0652  33.92815

# Row 6 (y=5225533 EMU)
Resize-Label "TextBox 23" 317.89406 409.03595 65.8874   33.92815
Resize-Label "TextBox 24" 468.7841  409.03595 98.9571   33.92815
Resize-Label "TextBox 25" 624.15485 409.03595 123.0652  33.92815
